$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: conversion message cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$nl = [char]10
$conversionText = "Conversión del día 💰" + $nl + `
"✅ Dólar paralelo: 68" + $nl + `
"" + $nl + `
"Binance" + $nl + `
"✅ 1000 Bs = 2.35 = 8870.59 pesos" + $nl + `
"✅ 8870.59 pesos = 2.33 = 953.76 Bs" + $nl + `
"" + $nl + `
"Promedio competencia" + $nl + `
"✅ Tasa pesos: 20" + $nl + `
"✅ Tasa Bs: 20" + $nl + `
"✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $conversionText

# --- Update the "tasas" sheet: N10, O10, N12, O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 425
$wsTasas.Range("O10").Value = 3770
$wsTasas.Range("N12").Value = 3804.9
$wsTasas.Range("O12").Value = 409.1
